$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-08 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-09 Saturday", 2) | Out-Null
$d.Content.Find.Execute("589÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "514÷8=", 2) | Out-Null
$d.Content.Find.Execute("647÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "277÷6=", 2) | Out-Null
$d.Content.Find.Execute("636÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "299÷2=", 2) | Out-Null
$d.Content.Find.Execute("310÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "522÷7=", 2) | Out-Null
$d.Content.Find.Execute("338÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "853÷5=", 2) | Out-Null
$d.Content.Find.Execute("692÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "786÷2=", 2) | Out-Null
$d.Content.Find.Execute("423÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "355÷5=", 2) | Out-Null
$d.Content.Find.Execute("259÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "200÷7=", 2) | Out-Null
$d.Content.Find.Execute("415÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "889÷5=", 2) | Out-Null
$d.Content.Find.Execute("369÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "285÷9=", 2) | Out-Null
$d.Content.Find.Execute("141÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "755÷2=", 2) | Out-Null
$d.Content.Find.Execute("597÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "647÷9=", 2) | Out-Null
$d.Content.Find.Execute("587÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "546÷8=", 2) | Out-Null
$d.Content.Find.Execute("703÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "476÷5=", 2) | Out-Null
$d.Content.Find.Execute("114÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "243÷7=", 2) | Out-Null
$d.Content.Find.Execute("299÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "736÷6=", 2) | Out-Null
$d.Content.Find.Execute("420÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "437÷6=", 2) | Out-Null
$d.Content.Find.Execute("544÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "661÷3=", 2) | Out-Null
$d.Content.Find.Execute("432÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "486÷8=", 2) | Out-Null
$d.Content.Find.Execute("474÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "138÷4=", 2) | Out-Null
$d.Content.Find.Execute("708÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "388÷8=", 2) | Out-Null
$d.Content.Find.Execute("914÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "648÷7=", 2) | Out-Null
$d.Content.Find.Execute("827÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "850÷8=", 2) | Out-Null
$d.Content.Find.Execute("453÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "329÷9=", 2) | Out-Null
$d.Content.Find.Execute("137÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "677÷2=", 2) | Out-Null
